$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The paragraph currently reads (one run, colour 000000):
#   "...ceste tourtelle plante ton animal, <del>aya</del> imitant..."
# We need to turn the literal text "ourtelle plante ton animal, " into
# three runs:
#   1) "ourtelle<comment>c_112v_0"   (same formatting as before: colour 000000)
#   2) "3"                           (no explicit colour - just rtl)
#   3) "</comment> plante ton animal, " (same formatting as before: colour 000000)
# i.e. insert a literal "<comment>c_112v_03</comment>" rendition-spec
# tag right after "ourtelle", with the digit "3" carrying slightly
# different run formatting than its neighbours.
# ------------------------------------------------------------------

$oldText = "ourtelle plante ton animal, "
$newText = "ourtelle<comment>c_112v_03</comment> plante ton animal, "

$content = $d.Content.Text
$start = $content.IndexOf($oldText)
$target = $d.Range($start, $start + $oldText.Length)
$target.Text = $newText

# Re-locate the digit "3" that sits inside the "<comment>...</comment>"
# tag we just inserted, so it can be given its own run.
$content = $d.Content.Text
$tagStart = $content.IndexOf("c_112v_03</comment>")
$digitPos = $tagStart + "c_112v_0".Length
$digitRun = $d.Range($digitPos, $digitPos + 1)

# Borrow the formatting (rPr with only <w:rtl w:val="0"/>, i.e. no
# explicit colour) from the existing "<del>aya</del>" run that already
# sits a little further along in this same paragraph - this is the
# only way, via the exposed object model, to produce a run whose rPr
# omits <w:color> altogether rather than writing w:color w:val="auto".
$content = $d.Content.Text
$donorStart = $content.IndexOf("<del>aya</del>") + "<del>".Length
$donorChar = $d.Range($donorStart, $donorStart + 1)

$digitRun.FormattedText = $donorChar.FormattedText
$digitRun.Text = "3"
